# Updated cryptos list on Tue Mar  5 15:00:16 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.726.11"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +3.81%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.796.94"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +7.39%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.03%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'419.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.46%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'137.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +4.51%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.784.65"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +7.13%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.649"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.05%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  +0.03%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.777"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.96%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +12.99%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.0000411"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +51.22%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'43.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -0.26%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'  +4.38%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'4.386.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +7.39%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  -0.57%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.821.61"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +8.75%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'20.55"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  +3.09%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'  +2.39%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'68.633.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.93%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'448.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.40%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'15.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +17.26%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'90.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.40%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'3.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -4.78%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'38.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +12.18%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'  -1.41%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.36%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'5.09"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +5.17%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +5.26%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'12.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  +1.48%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.81%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'7.11"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.12%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'41.76"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +6.44%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'58.01"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +1.10%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +0.13%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.0491"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.60%  "
$ws.Range("E38").Style = "Normal"

# Row 39/40 swap: ThetaToken <-> PEPE
$ws.Range("B39").Value = "'PEPE"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").Value = "'0.0₃0731"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.92%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("B40").Value = "'ThetaToken"
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").Value = "'3.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +29.89%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.67%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.997"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.18%  "
$ws.Range("E42").Style = "Normal"

# Row 43: LidoDAOToken -> ApeXProtocol
$ws.Range("B43").Value = "'ApeXProtocol"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'3.25"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  +27.21%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'27.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +28.06%  "
$ws.Range("E44").Style = "Normal"

# Row 45: Monero -> LidoDAOToken
$ws.Range("B45").Value = "'LidoDAOToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'3.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +4.15%  "
$ws.Range("E45").Style = "Normal"

# Row 46: ApeXProtocol -> Monero
$ws.Range("B46").Value = "'Monero"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'148.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.27%  "
$ws.Range("E46").Style = "Normal"

# Row 47: Stacks -> ARBITRUM
$ws.Range("B47").Value = "'ARBITRUM"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'2.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.00%  "
$ws.Range("E47").Style = "Normal"

# Row 48: ARBITRUM -> Stacks
$ws.Range("B48").Value = "'Stacks"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'2.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -5.07%  "
$ws.Range("E48").Style = "Normal"

# Row 49: NEARProtocol -> WEMIXToken
$ws.Range("B49").Value = "'WEMIXToken"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'2.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -5.43%  "
$ws.Range("E49").Style = "Normal"

# Row 50: WEMIXToken -> NEARProtocol
$ws.Range("B50").Value = "'NEARProtocol"
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = "'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").Value = "'4.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.83%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.303"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -2.53%  "
$ws.Range("E51").Style = "Normal"
